# Bump the "Förändrad" (Changed) date column from 2023-09-19 (45188)
# to 2023-09-20 (45189) for every data row on the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C262").Value = 45189
